# Regenerate save_data to use K instead of Strike#, recalculated std/mean, s_vals written.
# Only column G ("K") values change for rows 2-29 and 31-33 (rows 30, 34, 35 already correct).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 0
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 2
    19 = 1
    20 = 0
    21 = 0
    22 = 3
    23 = 0
    24 = 2
    25 = 1
    26 = 1
    27 = 0
    28 = 1
    29 = 1
    31 = 1
    32 = 2
    33 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
